$d = $word.ActiveDocument

# "Versi" + "on" -> single run "Version" (Find/Replace merges the runs
# spanned by the match into one, keeping the surrounding proofErr marks).
$rng1 = $d.Content
$rng1.Find.Execute("Version", $false, $false, $false, $false, $false, $true, 1, $false, "Version", 2)

# " 2" -> " 1." ; stop right before the trailing "." run so the
# bookmark (_GoBack) that sits between " 2" and "." is left untouched.
$rng2 = $d.Content
$rng2.Find.Execute(" 2", $false, $false, $false, $false, $false, $true, 1, $false, " 1.", 2)

# The original trailing "." run (now a redundant extra period after the
# bookmark) must be deleted so the paragraph reads "Version 1." again.
$p = $d.Paragraphs(1).Range
$trailingDot = $d.Range($p.End - 2, $p.End - 1)
$trailingDot.Text = ""
